$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Custom")

# Generalize AutoHotKey mappings: shift the F-key macro labels up by one
$ws.Range("B10").Value2 = "F14 (KeePassXC)"
$ws.Range("B11").Value2 = "F15 (Chat App)"

# Update the active selection on the Custom sheet to match the new state
$ws.Range("C12").Select()
